$wb = $excel.ActiveWorkbook

# "Subs" sheet: clear the countUp flag in F3 (Research row) and sync the selection
$wsSubs = $wb.Worksheets.Item("Subs")
$wsSubs.Range("F3").ClearContents()
$wsSubs.Activate()
$wsSubs.Range("E3").Select()
